# Insert a new data row at row 96 (pushing the existing rows 96..213 down to 97..214)
# and populate it with the new "Packham's Triumph / Primera" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("96:96").Insert()

$ws.Range("A96").Value = 7
$ws.Range("B96").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C96").Value = "Ñuble"
$ws.Range("D96").Value = 44781
$ws.Range("E96").Value = 16
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100104
$ws.Range("H96").Value = "Frutos de pepita"
$ws.Range("I96").Value = 100104005
$ws.Range("J96").Value = "Pera"
$ws.Range("K96").Value = "Packham's Triumph"
$ws.Range("L96").Value = "Primera"
$ws.Range("M96").Value = 120
$ws.Range("N96").Value = 9000
$ws.Range("O96").Value = 10000
$ws.Range("P96").Value = 9500
$ws.Range("Q96").Value = '$/caja 16 kilos empedrada'
$ws.Range("R96").Value = "Provincia de Curicó"
$ws.Range("S96").Value = 594
$ws.Range("T96").Value = 16
